$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
try {
  $ws.Rows.Item(581).Height = 60
  Write-Output "Height set ok"
} catch {
  Write-Output ("Error: " + $_.Exception.Message)
}
